$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content for rows 3, 4, 5 (Timeline slider + Symptom checker test cases) ---

# Row 3 - Timeline slider
$ws.Range("G3").Value = "We will test to ensure the images and content that are picked within the timeline slider display properly."
$ws.Range("H3").Value = "The timeline slider display in proper order without errors"
$ws.Range("I3").Value = "The images and contents do not advance through the sense, as they should. Display them out of the proper sequence without the designed interval; any of the images are not properly called from the database."

# Row 4 - Symptom checker / f&q generator
$ws.Range("G4").Value = "We will create a f&q generator on the page that helps to identify patient's health issue. The multiple choice form will consist of several radio buttons with the submissions posting back to the page using ajax to display aggregate result of overall submission."
$ws.Range("H4").Value = "The f&q generator accepts the users input and successfully passes it to the database which return the suitable health assessment."
$ws.Range("I4").Value = "The generator doesn't pass user input to database. It doesn't display the result properly."

# Row 5 - Content display
$ws.Range("G5").Value = "We will test to ensure the content are display properly"
$ws.Range("H5").Value = "The content display in proper order without errors"
$ws.Range("I5").Value = "The content do not advance thourh the sense as they should. Display them out of the properly sequence without the designed interval any of the images are not properly called from the database."

# --- Font styling: give the new cells a Cambria 12pt font, vertical-center + wrap ---
$ws.Range("I3:I5").Font.Name = "Cambria"
$ws.Range("I3:I5").Font.Size = 12
$ws.Range("I3:I5").VerticalAlignment = -4108
$ws.Range("I3:I5").WrapText = $true

$ws.Range("G4:H5").Font.Name = "Cambria"
$ws.Range("G4:H5").Font.Size = 12
$ws.Range("G4:H5").VerticalAlignment = -4108
$ws.Range("G4:H5").WrapText = $true

# G3/H3 keep the plain wrap-text style (same as column D entries)
$ws.Range("G3:H3").WrapText = $true

# --- Row heights for the new/updated rows ---
$ws.Rows.Item(3).RowHeight = 135
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 120
$ws.Rows.Item(6).RowHeight = 135

# --- View / selection changes ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("I5").Select()

$wb.Windows.Item(1).ScrollColumn = 2

Write-Host "done"
